# [FEATURE] Add recipe view to prototype
# Adds a new "Rezeptansicht Mock Up finished" task row to the Arbeitsmatrix
# sheet (row 123), pushes the trailing blank spacer rows down, and shifts
# the summary block (previously rows 125-128) down to rows 130-133,
# updating the formulas so they reference the new row numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# --- shift the summary block down by 5 rows (125-128 -> 130-133) -------
# Row 128 -> 133
$ws.Cells.Item(133, 2).Value = $ws.Cells.Item(128, 2).Value
$ws.Cells.Item(133, 3).Value = $ws.Cells.Item(128, 3).Value

# Row 127 -> 132
$ws.Cells.Item(132, 2).Value = $ws.Cells.Item(127, 2).Value
$ws.Cells.Item(132, 3).Formula = "=ROUNDUP(C130/30, 0)"

# Row 126 -> 131
$ws.Cells.Item(131, 4).Value = $ws.Cells.Item(126, 4).Value
$ws.Cells.Item(131, 5).Formula = "=135-E130"
$ws.Cells.Item(131, 6).Value = $ws.Cells.Item(126, 6).Value
$ws.Cells.Item(131, 7).Formula = "=315-G130"

# Row 125 -> 130
$ws.Cells.Item(130, 2).Value = $ws.Cells.Item(125, 2).Value
$ws.Cells.Item(130, 3).Formula = "=SUM(I:I)+SUM(H:H)"
$ws.Cells.Item(130, 4).Value = $ws.Cells.Item(125, 4).Value
$ws.Cells.Item(130, 5).Formula = "=SUM(H:H)"
$ws.Cells.Item(130, 6).Value = $ws.Cells.Item(125, 6).Value
$ws.Cells.Item(130, 7).Formula = "=SUM(I:I)"

# clear the old (now superseded) summary cells in rows 125-128 so they
# become the blank spacer rows seen in the target layout
$ws.Range("B125:G128").ClearContents()

# --- blank spacer rows 124-129 (D/F/G formatted, otherwise empty) ------
$ws.Range("D124:D129").Value = ""
$ws.Range("F124:F129").Value = ""
$ws.Range("G124:G129").Value = ""

# --- new task row 123: Rezeptansicht Mock Up finished ------------------
$ws.Cells.Item(123, 1).Value = 22
$ws.Cells.Item(123, 2).Value = $ws.Cells.Item(118, 2).Value
$ws.Cells.Item(123, 3).Value = $ws.Cells.Item(118, 3).Value
$ws.Cells.Item(123, 4).Value = $ws.Cells.Item(119, 4).Value
$ws.Cells.Item(123, 5).Value = "Rezeptansicht Mock Up finished"
$ws.Cells.Item(123, 6).Value = 44445
$ws.Cells.Item(123, 7).Value = 44481
$ws.Cells.Item(123, 9).Formula = "=ROUNDUP(((SUM(K123-J123)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(123, 10).Value = 0.38541666666666669
$ws.Cells.Item(123, 11).Value = 0.44791666666666669

# --- sheet view: scroll position + active selection ---------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 111
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F126").Select()

# --- data validation range grows from D115:D124 to D115:D129 -----------
$ws.Range("D2:D3,D13:D17,D22:D27,D35:D36,D41:D113,D115:D129").Validation.Delete()
$ws.Range("D2:D3,D13:D17,D22:D27,D35:D36,D41:D113,D115:D129").Validation.Add(3, 1, 1, "=$N$3:$N$6")
$ws.Range("D2:D3,D13:D17,D22:D27,D35:D36,D41:D113,D115:D129").Validation.ErrorTitle = "Prefix nicht unterstützt"
$ws.Range("D2:D3,D13:D17,D22:D27,D35:D36,D41:D113,D115:D129").Validation.ErrorMessage = "Es konnte kein korrekter Prefix ausgegeben werden`n"
$ws.Range("D2:D3,D13:D17,D22:D27,D35:D36,D41:D113,D115:D129").Validation.InputTitle = "Prefix"
$ws.Range("D2:D3,D13:D17,D22:D27,D35:D36,D41:D113,D115:D129").Validation.InputMessage = "Wählen Sie einen Prefix aus"

$wb.Save()
